$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the two trailing rows (24 "Requisitos:" and 25 the LOM3013 text) ---
# Their content is being folded into rows 22/23 below, so the sheet shrinks
# from A1:C25 to A1:C23.
$ws.Range("A24:C25").EntireRow.Delete()

# --- Row 10: Objetivos value is replaced by the professor line ---
$ws.Range("B10").Value = "5983729 - Fernando Vernilli Junior"
$ws.Range("C10").Value = "5983729 - Fernando Vernilli Junior"

# --- Row 13: now "Programa resumido:" with the activation date value ---
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "01/01/2020"
$ws.Range("C13").Value = "01/01/2020"
$ws.Rows.Item(13).RowHeight = 60

# --- Row 14: now just "Short syllabus:" (old professor values cleared) ---
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").ClearContents()
$ws.Range("C14").ClearContents()
$ws.Rows.Item(14).RowHeight = 60

# --- Row 15: now "Programa:" with the professor line ---
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "5983729 - Fernando Vernilli Junior"
$ws.Range("C15").Value = "5983729 - Fernando Vernilli Junior"
$ws.Rows.Item(15).RowHeight = 120

# --- Row 16: now "Syllabus:" ---
$ws.Range("A16").Value = "Syllabus:"
$ws.Rows.Item(16).RowHeight = 120

# --- Row 17: now "Avaliação:" only (old long Programa text removed) ---
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17").ClearContents()
$ws.Range("C17").ClearContents()
$ws.Rows.Item(17).AutoFit()

# --- Row 18: now "Método:" with the second professor line ---
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "1922320 - Sebastiao Ribeiro"
$ws.Range("C18").Value = "1922320 - Sebastiao Ribeiro"
$ws.Rows.Item(18).RowHeight = 60

# --- Row 19: now "Critério:" with the exam-method paragraph ---
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Serão realizadas duas provas escritas (P1 e P2), apresentações orais de trabalhos (T) e listas de exercícios (E)A nota final será calculada utilizando a equação: {[(P1 + P2 + T)/3] x 0,9} + E x 0,1"
$ws.Range("C19").Value = "Serão realizadas duas provas escritas (P1 e P2), apresentações orais de trabalhos (T) e listas de exercícios (E)A nota final será calculada utilizando a equação: {[(P1 + P2 + T)/3] x 0,9} + E x 0,1"
$ws.Rows.Item(19).RowHeight = 60

# --- Row 20: now "Norma de recuperação:" with the weighted-average criteria text ---
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "A nota final será a média ponderada das provas escritas (80% da nota final) e das listas de exercícios e relatórios (20% da nota final)."
$ws.Range("C20").Value = "A nota final será a média ponderada das provas escritas (80% da nota final) e das listas de exercícios e relatórios (20% da nota final)."
$ws.Rows.Item(20).RowHeight = 60

# --- Row 21: now "Bibliografia:" with the recovery-exam text ---
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "Para a recuperação será realizada uma prova (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez). NF = (MP + PR)/2. NF igual ou superior a 5 (cinco): aprovado. NF inferior a 5: reprovado."
$ws.Range("C21").Value = "Para a recuperação será realizada uma prova (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez). NF = (MP + PR)/2. NF igual ou superior a 5 (cinco): aprovado. NF inferior a 5: reprovado."
$ws.Rows.Item(21).RowHeight = 120

# --- Row 22: now just "Requisitos:" (old bibliography text cleared) ---
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()
$ws.Rows.Item(22).AutoFit()

# --- Row 23: now just the weak-requisite line in B/C (no label in A) ---
$ws.Range("A23").ClearContents()
$ws.Range("B23").Value = "LOM3013 -  Ciência dos Materiais  (Requisito fraco)`n"
$ws.Range("C23").Value = "LOM3013 -  Ciência dos Materiais  (Requisito fraco)`n"
$ws.Rows.Item(23).RowHeight = 30
